# Apply crypto price/volume updates scraped on Mon Jan  8 01:08:51 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.132.96"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "2.202.44"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'299.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.63%  "

$ws.Range("D6").Value = "'89.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.72%  "

$ws.Range("D7").Value = "'0.581"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.35%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.92%  "

$ws.Range("D10").Value = "'33.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.62%  "

$ws.Range("D11").Value = "'0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.20%  "

$ws.Range("D12").Value = "'6.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.56%  "

$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("D14").Value = "2.537.30"
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("D15").Value = "2.251.87"
$ws.Range("E15").Value = "  -4.51%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'13.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.783"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.41%  "

$ws.Range("D18").Value = "44.026.20"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").Value = "0.0₃0893"
$ws.Range("E19").Value = "  -7.47%  "

$ws.Range("D20").Value = "'5.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.18%  "

$ws.Range("D21").Value = "'11.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.55%  "

$ws.Range("D22").Value = "'64.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("D23").Value = "'231.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.71%  "

$ws.Range("D24").Value = "'2.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.70%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "'1.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.19%  "

$ws.Range("D27").Value = "'2.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.95%  "

$ws.Range("D28").Value = "'36.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.06%  "

$ws.Range("D29").Value = "'9.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.67%  "

$ws.Range("D30").Value = "'19.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.10%  "

$ws.Range("D31").Value = "'5.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.22%  "

$ws.Range("D32").Value = "'145.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.52%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0749"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.73%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.27%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.78%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.105"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.14%  "

$ws.Range("D38").Value = "'1.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.76%  "

$ws.Range("D39").Value = "'13.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.56%  "

$ws.Range("D40").Value = "'3.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.34%  "

$ws.Range("D41").Value = "'3.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.43%  "

$ws.Range("D42").Value = "'0.0284"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.03%  "

$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").Value = "1.752.31"
$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("D45").Value = "'1.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("B46").Value = "ordi"
$ws.Range("C46").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D46").Value = "'69.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'75.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.75%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.176"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.59%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'93.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.76%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.420.53"
$ws.Range("E50").Value = "  -2.17%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'7.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.10%  "

